# Adds a "Debug connector" (J3) line to the BOM, adds per-row "Notes" entries
# for several existing rows, tweaks one note's wording, and rebuilds the
# revision-history block (now with a dated "Rev | Date | Notes" table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Notes" (column J) entries for existing BOM rows.
# ---------------------------------------------------------------------
$ws.Range("J7").Value  = "RF connectors"
$ws.Range("J8").Value  = "Auxiliary (not placed!)"
$ws.Range("J9").Value  = "Main card connectors"
$ws.Range("J10").Value = "High wattage just in case"
$ws.Range("J12").Value = "-Z End Cap connector"

# ---------------------------------------------------------------------
# 2. New BOM row 13: debug connector (J3).
# ---------------------------------------------------------------------
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "J3"
$ws.Range("C13").Value = "P"
$ws.Range("C13").HorizontalAlignment = -4108   # xlCenter, matches other P/DNP cells
$ws.Range("D13").Value = "Harwin"
$ws.Range("E13").Value = "M50-3501042"
$ws.Range("F13").Value = '20 Position Connector Header Through Hole 0.050" (1.27mm)'
$ws.Range("G13").Value = "Digi-Key"
$ws.Range("H13").Value = "952-1386-ND"
$ws.Range("I13").Value = "N"
$ws.Range("J13").Value = "Debug connector"

# ---------------------------------------------------------------------
# 3. Rebuild the revision-history block, now one row lower (room for the
#    new row 13 pushed the whole block down), with an extra "Date" column.
# ---------------------------------------------------------------------

# Row 14 used to hold the "Revision versioing infro" caption -- clear it,
# it now lives on row 15.
$ws.Range("A14").ClearContents()

$ws.Range("A15").Value = "Revision versioing infro"

# Row 16 used to hold the lone "Rev" caption -- clear it, replaced by the
# 3-column header on row 17.
$ws.Range("A16").ClearContents()

$ws.Range("A17").Value = "Rev"
$ws.Range("B17").Value = "Date"
$ws.Range("C17").Value = "Notes"

# Row 18: first revision entry. No date recorded for it, so "--" is used
# as a placeholder in the (now date-formatted) Date column.
$ws.Range("A18").Value = "1.0r0"
$ws.Range("B18").Value = "--"
$ws.Range("C18").Value = "Initial BOM for 1U."

# Row 19: second revision entry, now mentions the debug connector too, and
# gets a real date.
$ws.Range("A19").Value = "1.1r0"
$ws.Range("B19").Value = 44173
$ws.Range("C19").Value = "Added -Z end cap connector, and debug connector."

# Apply the date number format across the whole prepared column (existing
# + blank rows reserved for future revisions), matching the author's
# selection/format sweep over B18:B36.
$ws.Range("B18:B36").NumberFormat = "yyyy-mm-dd"

# Re-assert the literal values after the bulk format (Date column's "--"
# must stay text, not get reinterpreted).
$ws.Range("B18").Value = "--"
$ws.Range("B19").Value = 44173

# Leave B20:B36 blank (formatted only), matching the saved selection.
$ws.Range("B18:B36").Select()
